$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 410.1
$ws.Range("I2").Value = 400.1111
$ws.Range("K2").Value = 400.1111
$ws.Range("M2").Value = -287.1111

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H43").Value = 37667.332
$ws.Range("I43").Value = 6500
$ws.Range("J43").Value = 100002
$ws.Range("K43").Value = 6500
$ws.Range("L43").Value = 100002
$ws.Range("M43").Value = -6431
$ws.Range("N43").Value = -100140

$ws.Range("H55").Value = 151.33333
$ws.Range("I55").Value = 139.09091
$ws.Range("J55").Value = 185
$ws.Range("K55").Value = 139.09091
$ws.Range("L55").Value = 185
$ws.Range("M55").Value = 74.90908999999999
$ws.Range("N55").Value = -613

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1382.7646
$ws.Range("I2").Value = 1024.5
$ws.Range("J2").Value = 2242.6
$ws.Range("K2").Value = 1024.5
$ws.Range("L2").Value = 2242.6
$ws.Range("M2").Value = -911.5
$ws.Range("N2").Value = -2468.6

$ws.Range("H44").Value = 28428.285
$ws.Range("J44").Value = 28428.285
$ws.Range("L44").Value = 28428.285
$ws.Range("N44").Value = -29404.285

$ws.Range("H51").Value = 21292.428
$ws.Range("I51").Value = 12000
$ws.Range("J51").Value = 25009.4
$ws.Range("K51").Value = 12000
$ws.Range("L51").Value = 25009.4
$ws.Range("M51").Value = -11244
$ws.Range("N51").Value = -26521.4

$ws.Range("H55").Value = 26735
$ws.Range("J55").Value = 26735
$ws.Range("L55").Value = 26735
$ws.Range("N55").Value = -27365

$ws.Range("H61").Value = 20835506
$ws.Range("I61").Value = 21741332
$ws.Range("K61").Value = 21741332
$ws.Range("M61").Value = -21741120

$ws.Range("H74").Value = 125001660
$ws.Range("I74").Value = 166667890
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 166667890
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -166667016
$ws.Range("N74").Value = -4748

$ws.Range("H77").Value = 125001660
$ws.Range("I77").Value = 166667890
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 833339450
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -833335082
$ws.Range("N77").Value = -23736

$ws.Range("H116").Value = 1382.7646
$ws.Range("I116").Value = 1024.5
$ws.Range("J116").Value = 2242.6
$ws.Range("K116").Value = 1024.5
$ws.Range("L116").Value = 2242.6
$ws.Range("M116").Value = 1269.5
$ws.Range("N116").Value = -6830.6

$ws.Range("H132").Value = 10001687
$ws.Range("I132").Value = 12501208
$ws.Range("J132").Value = 3602.4
$ws.Range("K132").Value = 37503624
$ws.Range("L132").Value = 10807.2
$ws.Range("M132").Value = -37501094
$ws.Range("N132").Value = -15867.2

$ws.Range("H136").Value = 20835506
$ws.Range("I136").Value = 21741332
$ws.Range("K136").Value = 65223996
$ws.Range("M136").Value = -65221446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1382.7646
$ws.Range("I3").Value = 1024.5
$ws.Range("J3").Value = 2242.6
$ws.Range("K3").Value = 1024.5
$ws.Range("L3").Value = 2242.6
$ws.Range("M3").Value = -910.5
$ws.Range("N3").Value = -2470.6

$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H56").Value = 46066
$ws.Range("J56").Value = 46066
$ws.Range("L56").Value = 46066
$ws.Range("N56").Value = -47544

$ws.Range("H86").Value = 18520896
$ws.Range("I86").Value = 2148.1428
$ws.Range("K86").Value = 2148.1428
$ws.Range("M86").Value = -1025.1428

$ws.Range("H89").Value = 18520896
$ws.Range("I89").Value = 2148.1428
$ws.Range("K89").Value = 10740.714
$ws.Range("M89").Value = -5124.714

$ws.Range("H94").Value = 900.55554
$ws.Range("I94").Value = 762.8570999999999
$ws.Range("J94").Value = 1382.5
$ws.Range("K94").Value = 762.8570999999999
$ws.Range("L94").Value = 1382.5
$ws.Range("M94").Value = -311.8570999999999
$ws.Range("N94").Value = -2284.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 27525.5
$ws.Range("J48").Value = 27525.5
$ws.Range("L48").Value = 27525.5
$ws.Range("N48").Value = -28477.5

$ws.Range("H132").Value = 41669080
$ws.Range("I132").Value = 62501828
$ws.Range("K132").Value = 187505484
$ws.Range("M132").Value = -187502954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 841.86957
$ws.Range("I5").Value = 311.4375
$ws.Range("J5").Value = 2054.2856
$ws.Range("K5").Value = 934.3125
$ws.Range("L5").Value = 6162.8568
$ws.Range("M5").Value = -822.3125
$ws.Range("N5").Value = -6386.8568

$ws.Range("H113").Value = 886.4595
$ws.Range("I113").Value = 547.0833
$ws.Range("K113").Value = 1641.2499
$ws.Range("M113").Value = 528.7501

$ws.Range("H117").Value = 583
$ws.Range("I117").Value = 250
$ws.Range("J117").Value = 916
$ws.Range("K117").Value = 750
$ws.Range("L117").Value = 2748
$ws.Range("M117").Value = 2692
$ws.Range("N117").Value = -9632

$ws.Range("H135").Value = 841.86957
$ws.Range("I135").Value = 311.4375
$ws.Range("J135").Value = 2054.2856
$ws.Range("K135").Value = 2802.9375
$ws.Range("L135").Value = 18488.5704
$ws.Range("M135").Value = -267.9375
$ws.Range("N135").Value = -23558.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6868.75
$ws.Range("I40").Value = 12000
$ws.Range("J40").Value = 4536.364
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 4536.364
$ws.Range("M40").Value = -11864
$ws.Range("N40").Value = -4808.364

$ws.Range("H60").Value = 25991
$ws.Range("J60").Value = 25991
$ws.Range("L60").Value = 25991
$ws.Range("N60").Value = -27009

$ws.Range("H68").Value = 1445.6666
$ws.Range("I68").Value = 1382
$ws.Range("J68").Value = 1545.7142
$ws.Range("K68").Value = 1382
$ws.Range("L68").Value = 1545.7142
$ws.Range("M68").Value = -633
$ws.Range("N68").Value = -3043.7142

$ws.Range("H71").Value = 1445.6666
$ws.Range("I71").Value = 1382
$ws.Range("J71").Value = 1545.7142
$ws.Range("K71").Value = 6910
$ws.Range("L71").Value = 7728.571
$ws.Range("M71").Value = -3166
$ws.Range("N71").Value = -15216.571

$ws.Range("H87").Value = 51496.332
$ws.Range("J87").Value = 51496.332
$ws.Range("L87").Value = 51496.332
$ws.Range("N87").Value = -53742.332

$ws.Range("H90").Value = 51496.332
$ws.Range("J90").Value = 51496.332
$ws.Range("L90").Value = 154488.996
$ws.Range("N90").Value = -165720.996

$ws.Range("H122").Value = 4751.5107
$ws.Range("I122").Value = 4746.852
$ws.Range("J122").Value = 4757.8
$ws.Range("K122").Value = 14240.556
$ws.Range("L122").Value = 14273.4
$ws.Range("M122").Value = -11790.556
$ws.Range("N122").Value = -19173.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2152.2856
$ws.Range("I100").Value = 3585.3333
$ws.Range("J100").Value = 1077.5
$ws.Range("K100").Value = 7170.6666
$ws.Range("L100").Value = 2155
$ws.Range("M100").Value = -3237
